# Replace the 25 "two-digit x two-digit" multiplication prompts in the
# practice-sheet table with the newly generated set of problems.
$d = $word.ActiveDocument

$d.Content.Find.Execute("23×89=", $true, $false, $false, $false, $false, $true, 1, $false, "57×20=", 2) | Out-Null
$d.Content.Find.Execute("87×20=", $true, $false, $false, $false, $false, $true, 1, $false, "62×58=", 2) | Out-Null
$d.Content.Find.Execute("87×59=", $true, $false, $false, $false, $false, $true, 1, $false, "44×34=", 2) | Out-Null
$d.Content.Find.Execute("55×29=", $true, $false, $false, $false, $false, $true, 1, $false, "97×68=", 2) | Out-Null
$d.Content.Find.Execute("65×47=", $true, $false, $false, $false, $false, $true, 1, $false, "72×56=", 2) | Out-Null
$d.Content.Find.Execute("16×57=", $true, $false, $false, $false, $false, $true, 1, $false, "67×31=", 2) | Out-Null
$d.Content.Find.Execute("15×77=", $true, $false, $false, $false, $false, $true, 1, $false, "22×32=", 2) | Out-Null
$d.Content.Find.Execute("55×69=", $true, $false, $false, $false, $false, $true, 1, $false, "47×94=", 2) | Out-Null
$d.Content.Find.Execute("27×18=", $true, $false, $false, $false, $false, $true, 1, $false, "35×23=", 2) | Out-Null
$d.Content.Find.Execute("31×80=", $true, $false, $false, $false, $false, $true, 1, $false, "81×64=", 2) | Out-Null
$d.Content.Find.Execute("78×97=", $true, $false, $false, $false, $false, $true, 1, $false, "12×52=", 2) | Out-Null
$d.Content.Find.Execute("27×19=", $true, $false, $false, $false, $false, $true, 1, $false, "95×81=", 2) | Out-Null
$d.Content.Find.Execute("47×35=", $true, $false, $false, $false, $false, $true, 1, $false, "17×95=", 2) | Out-Null
$d.Content.Find.Execute("54×70=", $true, $false, $false, $false, $false, $true, 1, $false, "61×35=", 2) | Out-Null
$d.Content.Find.Execute("91×60=", $true, $false, $false, $false, $false, $true, 1, $false, "17×69=", 2) | Out-Null
$d.Content.Find.Execute("90×29=", $true, $false, $false, $false, $false, $true, 1, $false, "26×97=", 2) | Out-Null
$d.Content.Find.Execute("63×96=", $true, $false, $false, $false, $false, $true, 1, $false, "70×77=", 2) | Out-Null
$d.Content.Find.Execute("26×27=", $true, $false, $false, $false, $false, $true, 1, $false, "19×98=", 2) | Out-Null
$d.Content.Find.Execute("66×50=", $true, $false, $false, $false, $false, $true, 1, $false, "98×63=", 2) | Out-Null
$d.Content.Find.Execute("62×81=", $true, $false, $false, $false, $false, $true, 1, $false, "58×30=", 2) | Out-Null
$d.Content.Find.Execute("93×36=", $true, $false, $false, $false, $false, $true, 1, $false, "82×52=", 2) | Out-Null
$d.Content.Find.Execute("34×55=", $true, $false, $false, $false, $false, $true, 1, $false, "39×98=", 2) | Out-Null
$d.Content.Find.Execute("14×13=", $true, $false, $false, $false, $false, $true, 1, $false, "63×13=", 2) | Out-Null
$d.Content.Find.Execute("74×67=", $true, $false, $false, $false, $false, $true, 1, $false, "90×84=", 2) | Out-Null
$d.Content.Find.Execute("42×50=", $true, $false, $false, $false, $false, $true, 1, $false, "66×35=", 2) | Out-Null
